# edit.ps1
# Applies the "poetry -> mathematics" content rewrite described by the diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a body fragment into the minimal "pkg:package" XML that
# Range.InsertXML expects, then insert it, replacing $rng's current content.
# ---------------------------------------------------------------------------
function Set-RangeBodyXml {
    param($rng, [string]$bodyXml)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

function R24 {
    param([string]$text, [string]$space = "")
    $spaceAttr = ""
    if ($space -eq "preserve") { $spaceAttr = ' xml:space="preserve"' }
    return '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t' + $spaceAttr + '>' + $text + '</w:t></w:r>'
}

# =============================================================================
# Paragraph 5 : the long body paragraph
# =============================================================================
$p5 = $d.Paragraphs.Item(5)
$body5 = '<w:p>'
$body5 += R24 'In this study of numbers, patterns, and pure logic, mathematics provides an analytical lens through which we comprehend the natural world and human constructs'
$body5 += R24 '.'
$body5 += R24 ' Through concepts like infinity, fractals, and the golden ratio, mathematics reveals hidden configurations and structures behind the facade of everyday life, inspiring awe and a sense of profound wonder' 'preserve'
$body5 += R24 '.'
$body5 += R24 ' From the harmonious ratios in music to the logarithmic spirals in nature, mathematics manifests itself as an omnipresent language that underlies the fabric of reality' 'preserve'
$body5 += R24 '.'
$body5 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r>'
$body5 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Mathematics offers a rich framework for solving problems intuitively and cultivating creativity</w:t></w:r>'
$body5 += R24 '.'
$body5 += R24 ' Algebra offers a powerful toolset for understanding patterns, capturing relationships between variables, and constructing abstract models that empower us to predict outcomes' 'preserve'
$body5 += R24 '.'
$body5 += R24 ' Calculus, with its derivatives and integrals, unlocks the dynamics of change and growth, providing a mathematical microscope into the behavior of systems over time' 'preserve'
$body5 += R24 '.'
$body5 += R24 ' The exploration of mathematical truths, theorems, and axioms, often stemming from seemingly simple postulates, reveals the underlying coherence and order inherent in our universe' 'preserve'
$body5 += R24 '.'
$body5 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r>'
$body5 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Moreover, mathematics transcends cultural boundaries and historical epochs</w:t></w:r>'
$body5 += R24 '.'
$body5 += R24 " Archimedes' treatise on floating bodies still holds true today; Euclid's geometry continues to furnish insights into shapes and angles; and Newton's formulation of gravity remains a cornerstone of modern physics" 'preserve'
$body5 += R24 '.'
$body5 += R24 ' It is a testament to the enduring power and universality of mathematics as a human endeavor, connecting individuals across time, culture, and geography in a shared intellectual dialogue' 'preserve'
$body5 += R24 '.'
$body5 += R24 ' By venturing into this realm of abstract beauty and pure thought, mathematics opens doors to intricate patterns, complex structures, and the profound satisfaction of intellectual discovery' 'preserve'
$body5 += R24 '.'
$body5 += '</w:p>'

Set-RangeBodyXml $p5.Range $body5

# =============================================================================
# Paragraph 7 : the "Summary" body paragraph (+ one new trailing empty para)
# =============================================================================
$p7 = $d.Paragraphs.Item(7)
$body7 = '<w:p>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>In this essay, we have explored the realm of mathematics, unraveling its complexities through concepts like infinity, fractals, and the golden ratio</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> We have delved into the analytical nature of mathematics and acknowledged its omnipresence in understanding the world through the harmonious ratios in music and logarithmic spirals in nature</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> The power of algebra, calculus, and </w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>the coherence of theorems and axioms were illuminated</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> We also recognized the universality and enduring nature of mathematics, transcending boundaries and connecting individuals in a shared intellectual dialogue</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Understanding mathematics is not merely about mastering equations and formulas but appreciating its inherent beauty and its potential to unravel the mysteries of the universe</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Through numbers, patterns, and logic, we gain access to a deeper understanding of the world around us, unlocking the gates of intellectual discovery and inspiring us with a sense of awe and wonder</w:t></w:r>'
$body7 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>'
$body7 += '</w:p>'
$body7 += '<w:p/>'

Set-RangeBodyXml $p7.Range $body7

# =============================================================================
# Title, subtitle, and email line (simple text swaps)
# =============================================================================
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Echoes of Ancient Lines in Modern Rhymes", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Exploring the Realm of Mathematics: A Journey through Numbers, Patterns, and Logic", 2) | Out-Null

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("Emily Dickinson", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Emily Anderson", 2) | Out-Null

# Email paragraph: "emilydickinson@poetry" + "." + "com"
#              ->  "emily" + "." + "anderson @ validweb" + "." + "com"
$p3 = $d.Paragraphs.Item(3)
$body3 = '<w:p><w:pPr><w:jc w:val="center"/></w:pPr>'
$body3 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>emily</w:t></w:r>'
$body3 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r>'
$body3 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>anderson @ validweb</w:t></w:r>'
$body3 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r>'
$body3 += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>com</w:t></w:r>'
$body3 += '</w:p>'
Set-RangeBodyXml $p3.Range $body3

# =============================================================================
# Font fix: TimesNewToman -> Times New Roman, everywhere
# =============================================================================
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}

Write-Output "done"
